# Update DKW1 electrolysis capacity (H77: 1500 -> 3500)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capacity")

$ws.Range("H77").Value = 3500

# Re-apply the AutoFilter on the Capacity sheet:
#  - extend the filtered range to the full data (now through row 172)
#  - drop the old Node (column A) filter for ES00/FI00
#  - filter Generator_ID (column C) on Electrolysis, Hydrogen processor,
#    Hydrogen storage dimensioner
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$filterRange = $ws.Range("A1:J172")
$filterRange.AutoFilter(3, @("Electrolysis", "Hydrogen processor", "Hydrogen storage dimensioner"), 7)

# Keep the hidden _xlnm._FilterDatabase defined name in sync with the new
# filter range (Excel normally does this automatically when the filter is
# applied interactively).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Capacity!_FilterDatabase") {
        $n.RefersTo = "=Capacity!`$A`$1:`$J`$172"
    }
}

# Move the active selection on the Capacity sheet
$ws.Range("B177").Select()
